$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added to the data set. They are
# inserted at rows 236-237, which pushes the previous rows 236-251 down
# to rows 238-253 (dimension grows from A1:T251 to A1:T253).
$ws.Rows("236:237").Insert()

# New record 1 (row 236): Especial quality, Provincia de Curicó origin.
$ws.Cells.Item(236, 1).Value = 6
$ws.Cells.Item(236, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(236, 3).Value = "Metropolitana"
$ws.Cells.Item(236, 4).Value = 44931
$ws.Cells.Item(236, 5).Value = 13
$ws.Cells.Item(236, 6).Value = "Fruta"
$ws.Cells.Item(236, 7).Value = 100101
$ws.Cells.Item(236, 8).Value = "Berries"
$ws.Cells.Item(236, 9).Value = 100101004
$ws.Cells.Item(236, 10).Value = "Frambuesa"
$ws.Cells.Item(236, 11).Value = "Sin especificar"
$ws.Cells.Item(236, 12).Value = "Especial"
$ws.Cells.Item(236, 13).Value = 300
$ws.Cells.Item(236, 14).Value = 8000
$ws.Cells.Item(236, 15).Value = 8000
$ws.Cells.Item(236, 16).Value = 8000
$ws.Cells.Item(236, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(236, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(236, 19).Value = 4000
$ws.Cells.Item(236, 20).Value = 2

# New record 2 (row 237): Especial quality, Región del Maule origin.
$ws.Cells.Item(237, 1).Value = 6
$ws.Cells.Item(237, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(237, 3).Value = "Metropolitana"
$ws.Cells.Item(237, 4).Value = 44931
$ws.Cells.Item(237, 5).Value = 13
$ws.Cells.Item(237, 6).Value = "Fruta"
$ws.Cells.Item(237, 7).Value = 100101
$ws.Cells.Item(237, 8).Value = "Berries"
$ws.Cells.Item(237, 9).Value = 100101004
$ws.Cells.Item(237, 10).Value = "Frambuesa"
$ws.Cells.Item(237, 11).Value = "Sin especificar"
$ws.Cells.Item(237, 12).Value = "Especial"
$ws.Cells.Item(237, 13).Value = 250
$ws.Cells.Item(237, 14).Value = 8000
$ws.Cells.Item(237, 15).Value = 8000
$ws.Cells.Item(237, 16).Value = 8000
$ws.Cells.Item(237, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(237, 18).Value = "Región del Maule"
$ws.Cells.Item(237, 19).Value = 4000
$ws.Cells.Item(237, 20).Value = 2
